# Reorders the comma-separated list of "Recorded By" values in column G
# so that any "System"/"system" tokens appear first (preserving their
# relative order), followed by the remaining tokens sorted alphabetically.
# This matches the canonical re-ordering performed by the upstream sync job.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Column G is the 7th column ("Recorded By")
$col = 7

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Text

    if ($null -eq $val) { continue }
    if (-not ($val -is [string])) { continue }
    if ($val -notmatch ",") { continue }

    $parts = $val -split "," | ForEach-Object { $_.Trim() }

    $systemParts = @()
    $otherParts = @()
    foreach ($p in $parts) {
        if ($p.ToLower() -eq "system") {
            $systemParts += $p
        } else {
            $otherParts += $p
        }
    }

    $otherSorted = $otherParts | Sort-Object
    $newParts = @()
    $newParts += $systemParts
    $newParts += $otherSorted

    $newVal = [string]::Join(", ", $newParts)

    if ($newVal -ne $val) {
        $cell.Value = $newVal
    }
}
